{"js": "// Fix the duplicated word \"Sternbild Sternbild\" -> \"Sternbild\" in the\n// German GaN 2022 Hercules campaign-dates sentence. The sentence occurs\n// multiple times throughout the document, so find every occurrence and\n// replace each one individually (format-preserving, single run per hit).\nconst searchResults = context.document.body.search(\"Sternbild Sternbild Herkules\", {\n  matchCase: true,\n  matchWholeWord: false\n});\nsearchResults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < searchResults.items.length; i++) {\n  searchResults.items[i].insertText(\"Sternbild Herkules\", \"Replace\");\n}\nawait context.sync();\n", "ps1": "# Fix the duplicated word \"Sternbild Sternbild\" -> \"Sternbild\" in the\n# German GaN 2022 Hercules campaign-dates sentence. The sentence is\n# repeated several times in the document body, so run Find/Replace over\n# the whole story range with Replace:=wdReplaceAll (2) to catch every hit.\n$d = $word.ActiveDocument\n$range = $d.Content\n$find = $range.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"Sternbild Sternbild Herkules\"\n$find.Replacement.Text = \"Sternbild Herkules\"\n$find.Forward = $true\n$find.Wrap = 1\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\n$find.Execute(\n    $find.Text,\n    $find.MatchCase,\n    $find.MatchWholeWord,\n    $find.MatchWildcards,\n    $false,\n    $false,\n    $find.Forward,\n    $find.Wrap,\n    $false,\n    $find.Replacement.Text,\n    2\n) | Out-Null\n"}
